$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E17: "Pendiente ADM" -> "ICD30759585"
$ws.Range("E17").Value = "ICD30759585"

# Remove row 55 (Caso 6572 / MEXICO 2639 / Almagro) - subsequent rows shift up
$ws.Rows(55).Delete()
